$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rows 96-99 previously belonged to "Seção" 4; they now belong to
#    Seção 3 (merged with the preceding block of rows).
# ------------------------------------------------------------------
$ws.Range("B96").Value = 3
$ws.Range("B97").Value = 3
$ws.Range("B98").Value = 3
$ws.Range("B99").Value = 3

# ------------------------------------------------------------------
# 2) Grow the worksheet table ("Tabela1") so it covers the two new
#    rows that are about to be appended.
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B1:G101"))

# ------------------------------------------------------------------
# 3) Add the two new data rows (100 and 101) describing the new
#    "Aula" about validating a non-duplicated e-mail on Cliente
#    update. Formatting (font/wrap/row height) is copied from the
#    most similar existing rows (98 and 97 respectively) so the new
#    rows render consistently with the rest of the table.
# ------------------------------------------------------------------
$ws.Range("B98:G98").Copy($ws.Range("B100:G100"))
$ws.Range("B97:G97").Copy($ws.Range("B101:G101"))

$ws.Range("B100").Value = 3
$ws.Range("C100").Value = "Operações de CRUD e Casos de Uso"
$ws.Range("D100").Value = 47

$ws.Range("B101").Value = 3
$ws.Range("C101").Value = "Operações de CRUD e Casos de Uso"
$ws.Range("D101").Value = 47

$ws.Range("E100").Value = " Validação customizada: email não repetido na atualização de Cliente"
$ws.Range("E101").Value = " Validação customizada: email não repetido na atualização de Cliente"

# Write F101 before F100 so the shared-string table gets the two new
# strings interned in the same order as the authored workbook.
$ws.Range("F101").Value = "5:13 - estrutura Map<String,String> - responsável por captar parametros recebidos pela URI"
$ws.Range("F100").Value = "4:50 - HttpServletRequest - tem função que permite obter parametros passados pela URI"

$ws.Range("G100").Value = "`n`n"
$ws.Range("G101").Value = ""

$ws.Rows.Item(100).RowHeight = 45
$ws.Rows.Item(101).RowHeight = 45

# ------------------------------------------------------------------
# 4) Update the sheet's last active selection to match the cell the
#    author was last editing.
# ------------------------------------------------------------------
$ws.Range("E100").Select() | Out-Null
